$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the account-summary totals (VALOR MORA, Cant. Trabajadores, Cant. Periodos) ---
$ws.Cells.Item(11, 5).Value2 = 122887   # E11: VALOR MORA total
$ws.Cells.Item(13, 3).Value2 = 2        # C13: Cant. Trabajadores
$ws.Cells.Item(13, 6).Value2 = 3        # F13: Cant. Periodos

# --- Replace the worker detail table ---
# The first three detail rows (CARMEN CECILIA CASTRO JIMENEZ, and LUISA
# FERNANDA VILLA JULIO x2) are removed entirely. Deleting them shifts the
# remaining two workers (JOANA MARCELA PEREZ and SUGEY DEL CARMEN MATURANA
# ROSENSTAND, previously in rows 19-21) up into rows 16-18, carrying their
# existing formatting (including the bottom-border style on the final row).
$ws.Rows("16:18").Delete()

# Row 16: JOANA MARCELA PEREZ, periodo 1901
$ws.Cells.Item(16, 3).Value2 = "53124729"
$ws.Cells.Item(16, 4).Value2 = "JOANA MARCELA PEREZ"
$ws.Cells.Item(16, 5).Value2 = "1901"
$ws.Cells.Item(16, 6).Value2 = 41377
$ws.Cells.Item(16, 7).Value2 = 1075809

# Row 17: JOANA MARCELA PEREZ, periodo 1902
$ws.Cells.Item(17, 3).Value2 = "53124729"
$ws.Cells.Item(17, 4).Value2 = "JOANA MARCELA PEREZ"
$ws.Cells.Item(17, 5).Value2 = "1902"
$ws.Cells.Item(17, 6).Value2 = 41377
$ws.Cells.Item(17, 7).Value2 = 1075809

# Row 18 (last data row, keeps the bottom-border style): SUGEY DEL CARMEN
# MATURANA ROSENSTAND, periodo 2206
$ws.Cells.Item(18, 3).Value2 = "33101883"
$ws.Cells.Item(18, 4).Value2 = "SUGEY DEL CARMEN MATURANA ROSENSTAND"
$ws.Cells.Item(18, 5).Value2 = "2206"
$ws.Cells.Item(18, 6).Value2 = 40133
$ws.Cells.Item(18, 7).Value2 = 6934400
